$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final roster (row 2 .. row 18): Name, Position, Team
$players = @(
    @("Immanuel Quickley", "PG,SG", "Toronto Raptors"),
    @("Desmond Bane", "SG,SF", "Memphis Grizzlies"),
    @("Kawhi Leonard", "SG,SF,PF", "LA Clippers"),
    @("Guerschon Yabusele", "PF,C", "Philadelphia 76ers"),
    @("D'Angelo Russell", "PG", "Brooklyn Nets"),
    @("Devin Booker", "PG,SG", "Phoenix Suns"),
    @("Trae Young", "PG", "Atlanta Hawks"),
    @("Jalen Brunson", "PG", "New York Knicks"),
    @("LeBron James", "SF,PF", "Los Angeles Lakers"),
    @("Christian Braun", "SG,SF", "Denver Nuggets"),
    @("Myles Turner", "C", "Indiana Pacers"),
    @("Walker Kessler", "C", "Utah Jazz"),
    @("Scoot Henderson", "PG", "Portland Trail Blazers"),
    @("Jalen Williams", "SG,SF,PF,C", "Oklahoma City Thunder"),
    @("Norman Powell", "SG,SF", "LA Clippers"),
    @("Brandon Ingram", "SG,SF,PF", "New Orleans Pelicans"),
    @("Jimmy Butler", "SF,PF", "Miami Heat")
)

$row = 2
foreach ($p in $players) {
    $ws.Cells.Item($row, 1).Value = $p[0]
    $ws.Cells.Item($row, 2).Value = $p[1]
    $ws.Cells.Item($row, 3).Value = $p[2]
    $row = $row + 1
}
